# Windows Python Http Server.docx edit
# Commit: "revise message article, fix an error on cookie, add a link"
#
# Content change: append "Under Windows 7" after the existing instruction
# about opening a command window, turning:
#   Shift + Right Mouse Click > Select "Open command window here"
# into:
#   Shift + Right Mouse Click > Select "Open command window here" Under Windows 7

$d = $word.ActiveDocument

$openQuote  = [char]8220
$closeQuote = [char]8221

$oldText = "Shift + Right Mouse Click > Select " + $openQuote + "Open command window here" + $closeQuote
$newText = "Shift + Right Mouse Click > Select " + $openQuote + "Open command window here" + $closeQuote + " Under Windows 7"

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2)
